$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3, 1, 0.4767983333333334, 1.430395, 0.003723890400117776, 0.003723890400117776, 3, 1, 32.63563666666667, 97.90691000000001, 0.9900101876891448, 0.9900101876891446, 15.56061716993889, 140.04555452945, 0.003686689433954404, 0.003686689433954404)
    3 = @(3, 1, 0.4767983333333334, 1.430395, 0.003723890400117776, 0.003723890400117776, 1, 0.3333333333333333, 0.05920466666666666, 0.177614, 0.001795988347259859, 0.001795988347259859, 0.02822868639222222, 0.25405817753, 0.000006688063765084379, 0.000006688063765084378)
    4 = @(3, 1, 0.4767983333333334, 1.430395, 0.003723890400117776, 0.003723890400117776, 3, 1, 0.270109, 0.810327, 0.008193823963595435, 0.008193823963595434, 0.1287875210183334, 1.159087689165, 0.00003051290239828803, 0.00003051290239828802)
    5 = @(3, 1, 108.7065656666667, 326.119697, 0.8490200321922391, 0.8490200321922391, 3, 1, 32.63563666666667, 97.90691000000001, 0.9900101876891448, 0.9900101876891446, 3547.707980378475, 31929.37182340627, 0.8405384814224824, 0.8405384814224822)
    6 = @(3, 1, 108.7065656666667, 326.119697, 0.8490200321922391, 0.8490200321922391, 1, 0.3333333333333333, 0.05920466666666666, 0.177614, 0.001795988347259859, 0.001795988347259859, 6.435935984773111, 57.923423862958, 0.001524830084407452, 0.001524830084407451)
    7 = @(3, 1, 108.7065656666667, 326.119697, 0.8490200321922391, 0.8490200321922391, 3, 1, 0.270109, 0.810327, 0.008193823963595435, 0.008193823963595434, 29.36262174565766, 264.263595710919, 0.006956720685349337, 0.006956720685349335)
    8 = @(3, 1, 18.85432833333333, 56.562985, 0.1472560774076432, 0.1472560774076432, 3, 1, 32.63563666666667, 97.90691000000001, 0.9900101876891448, 0.9900101876891446, 615.3230090807056, 5537.90708172635, 0.145785016832708, 0.145785016832708)
    9 = @(3, 1, 18.85432833333333, 56.562985, 0.1472560774076432, 0.1472560774076432, 1, 0.3333333333333333, 0.05920466666666666, 0.177614, 0.001795988347259859, 0.001795988347259859, 1.116264224198889, 10.04637801779, 0.0002644701990873228, 0.0002644701990873229)
    10 = @(3, 1, 18.85432833333333, 56.562985, 0.1472560774076432, 0.1472560774076432, 3, 1, 0.270109, 0.810327, 0.008193823963595435, 0.008193823963595434, 5.092723771788332, 45.834513946095, 0.001206590375847811, 0.001206590375847811)
}

foreach ($r in $data.Keys) {
    $rowVals = $data[$r]
    for ($i = 0; $i -lt $rowVals.Count; $i++) {
        $ws.Cells.Item([int]$r, 5 + $i).Value = $rowVals[$i]
    }
}

Write-Host "Updated rows 2-10, columns E:T"
